$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status of "FoodOrderingSystem" (row 4) from "Ongoing" to "Done"
$ws.Range("B4").Value = "Done"

# Update the selected cell/active cell shown in the saved view
$ws.Range("G6").Select()
